$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 currently empty; fill it in following the same pattern as row 10
# (date / request text / date-received / (blank) / observation text).
# C11 must pick up the "highlighted" date style used in column C (style of C10),
# so copy that cell's formatting over before setting values.
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A11").Value2 = 44578
$ws.Range("B11").Value2 = "Que se vea el campo Novedades en las Obras de Tasa"
$ws.Range("C11").Value2 = 44578
$ws.Range("E11").Value2 = " Cambio en API"

# Update the saved selection to A12, matching the new cursor position.
$ws.Range("A12").Select()
